# Scheduled runner update: refresh market-price-derived Leve profit figures
# across the per-job Leve Profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR).
# Values below are the latest currentAveragePrice-derived numbers for the
# affected rows; downstream profit columns (M/N) are updated to match.

$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets("ALC")
$ws.Range("H64").Value = 3790
$ws.Range("J64").Value = 3980
$ws.Range("L64").Value = 3980
$ws.Range("N64").Value = -4476
$ws.Range("H67").Value = 3790
$ws.Range("J67").Value = 3980
$ws.Range("L67").Value = 3980
$ws.Range("N67").Value = -5696
$ws.Range("H76").Value = 4899.8887
$ws.Range("I76").Value = 6033.3335
$ws.Range("J76").Value = 4333.1665
$ws.Range("K76").Value = 6033.3335
$ws.Range("L76").Value = 4333.1665
$ws.Range("M76").Value = -5718.3335
$ws.Range("N76").Value = -4963.1665
$ws.Range("H79").Value = 4899.8887
$ws.Range("I79").Value = 6033.3335
$ws.Range("J79").Value = 4333.1665
$ws.Range("K79").Value = 6033.3335
$ws.Range("L79").Value = 4333.1665
$ws.Range("M79").Value = -4941.3335
$ws.Range("N79").Value = -6517.1665
$ws.Range("H100").Value = 2125
$ws.Range("I100").Value = 1200
$ws.Range("K100").Value = 1200
$ws.Range("M100").Value = -659
$ws.Range("H112").Value = 2796.1667
$ws.Range("J112").Value = 3008.1875
$ws.Range("L112").Value = 9024.5625
$ws.Range("N112").Value = -11240.5625
$ws.Range("H129").Value = 826.9394
$ws.Range("I129").Value = 541.6667
$ws.Range("J129").Value = 890.3333
$ws.Range("K129").Value = 1625.0001
$ws.Range("L129").Value = 2670.9999
$ws.Range("M129").Value = 3374.9999
$ws.Range("N129").Value = -12670.9999
$ws.Range("H137").Value = 1420.262
$ws.Range("I137").Value = 1139.6111
$ws.Range("J137").Value = 1630.75
$ws.Range("K137").Value = 3418.8333
$ws.Range("L137").Value = 4892.25
$ws.Range("M137").Value = -868.8333000000002
$ws.Range("N137").Value = -9992.25
$ws.Range("H138").Value = 1837.3368
$ws.Range("I138").Value = 1250.3529
$ws.Range("J138").Value = 1965.2693
$ws.Range("K138").Value = 3751.0587
$ws.Range("L138").Value = 5895.8079
$ws.Range("M138").Value = 1388.9413
$ws.Range("N138").Value = -16175.8079
$ws = $wb.Sheets("ARM")
$ws.Range("H32").Value = 4857.354
$ws.Range("I32").Value = 4673.0483
$ws.Range("K32").Value = 4673.0483
$ws.Range("M32").Value = -4386.0483
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H61").Value = 100001400
$ws.Range("I61").Value = 142858140
$ws.Range("J61").Value = 2338
$ws.Range("K61").Value = 142858140
$ws.Range("L61").Value = 2338
$ws.Range("M61").Value = -142857928
$ws.Range("N61").Value = -2762
$ws.Range("H74").Value = 2959.9
$ws.Range("I74").Value = 1942.7142
$ws.Range("J74").Value = 5333.3335
$ws.Range("K74").Value = 1942.7142
$ws.Range("L74").Value = 5333.3335
$ws.Range("M74").Value = -1068.7142
$ws.Range("N74").Value = -7081.3335
$ws.Range("H77").Value = 2959.9
$ws.Range("I77").Value = 1942.7142
$ws.Range("J77").Value = 5333.3335
$ws.Range("K77").Value = 9713.571
$ws.Range("L77").Value = 26666.6675
$ws.Range("M77").Value = -5345.571
$ws.Range("N77").Value = -35402.6675
$ws.Range("H97").Value = 424.5
$ws.Range("I97").Value = 409.4
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 409.4
$ws.Range("L97").Value = 500
$ws.Range("M97").Value = 86.60000000000002
$ws.Range("N97").Value = -1492
$ws.Range("H132").Value = 2938.739
$ws.Range("I132").Value = 2504.3333
$ws.Range("K132").Value = 7512.999899999999
$ws.Range("M132").Value = -4982.999899999999
$ws.Range("H136").Value = 100001400
$ws.Range("I136").Value = 142858140
$ws.Range("J136").Value = 2338
$ws.Range("K136").Value = 428574420
$ws.Range("L136").Value = 7014
$ws.Range("M136").Value = -428571870
$ws.Range("N136").Value = -12114
$ws = $wb.Sheets("BSM")
$ws.Range("H105").Value = 53153308
$ws.Range("I105").Value = 56106164
$ws.Range("K105").Value = 56106164
$ws.Range("M105").Value = -56104417
$ws.Range("H134").Value = 1665.8
$ws.Range("I134").Value = 1239.75
$ws.Range("J134").Value = 1949.8334
$ws.Range("K134").Value = 3719.25
$ws.Range("L134").Value = 5849.5002
$ws.Range("M134").Value = -1184.25
$ws.Range("N134").Value = -10919.5002
$ws = $wb.Sheets("CRP")
$ws.Range("H38").Value = 2000
$ws.Range("J38").Value = 2000
$ws.Range("L38").Value = 2000
$ws.Range("N38").Value = -2754
$ws.Range("H46").Value = 2000
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2422
$ws.Range("H134").Value = 17859242
$ws.Range("I134").Value = 2095.1365
$ws.Range("J134").Value = 83335450
$ws.Range("K134").Value = 6285.4095
$ws.Range("L134").Value = 250006350
$ws.Range("M134").Value = -3750.4095
$ws.Range("N134").Value = -250011420
$ws = $wb.Sheets("CUL")
$ws.Range("H131").Value = 11237037
$ws.Range("J131").Value = 1133.2048
$ws.Range("L131").Value = 3399.6144
$ws.Range("N131").Value = -13479.6144
$ws = $wb.Sheets("GSM")
$ws.Range("H70").Value = 21432212
$ws.Range("I70").Value = 22730922
$ws.Range("J70").Value = 20003630
$ws.Range("K70").Value = 22730922
$ws.Range("L70").Value = 20003630
$ws.Range("M70").Value = -22730652
$ws.Range("N70").Value = -20004170
$ws.Range("H73").Value = 21432212
$ws.Range("I73").Value = 22730922
$ws.Range("J73").Value = 20003630
$ws.Range("K73").Value = 22730922
$ws.Range("L73").Value = 20003630
$ws.Range("M73").Value = -22729986
$ws.Range("N73").Value = -20005502
$ws.Range("H132").Value = 3412.9583
$ws.Range("I132").Value = 3298.0667
$ws.Range("K132").Value = 9894.2001
$ws.Range("M132").Value = -7364.2001
$ws = $wb.Sheets("WVR")
$ws.Range("H122").Value = 17858650
$ws.Range("I122").Value = 22728918
$ws.Range("K122").Value = 68186754
$ws.Range("M122").Value = -68184304
$ws.Range("H132").Value = 1246.6578
$ws.Range("I132").Value = 1061.5483
$ws.Range("J132").Value = 2066.4285
$ws.Range("K132").Value = 3184.6449
$ws.Range("L132").Value = 6199.2855
$ws.Range("M132").Value = -654.6448999999998
$ws.Range("N132").Value = -11259.2855
$ws.Range("H136").Value = 1093.4572
$ws.Range("I136").Value = 1015.32
$ws.Range("J136").Value = 1288.8
$ws.Range("K136").Value = 3045.96
$ws.Range("L136").Value = 3866.4
$ws.Range("M136").Value = -495.96
$ws.Range("N136").Value = -8966.4
